# Apply the two changes described by the commit:
#  1. Slide 6's table switches to a different (built-in) table style.
#  2. The presentation's theme colour scheme (the one actually driving the
#     slide master / slides) is changed from the custom "Integral" palette
#     back to the default "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{007C2D55-EDA1-4DF9-8463-9E5306CC092A}")
}

# --- 2. Swap the theme colour scheme back to the default Office palette --
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Length; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $themeColors.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}
